# Insert a new data row at row 20 (pushing existing rows 20-94 down to 21-95)
# and populate the newly inserted row with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 20 - this shifts rows 20:94 -> 21:95
$ws.Rows.Item(20).Insert()

# Populate the new row 20 with the new record's data
$ws.Range("A20").Value = 11
$ws.Range("B20").Value = 'Vega Monumental Concepción'
$ws.Range("C20").Value = 'Bíobío'
$ws.Range("D20").Value = 44624
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = 'Fruta'
$ws.Range("G20").Value = 100101
$ws.Range("H20").Value = 'Berries'
$ws.Range("I20").Value = 100101001
$ws.Range("J20").Value = 'Arándano (blue)'
$ws.Range("K20").Value = 'Sin especificar'
$ws.Range("L20").Value = 'Primera'
$ws.Range("M20").Value = 150
$ws.Range("N20").Value = 2500
$ws.Range("O20").Value = 2500
$ws.Range("P20").Value = 2500
$ws.Range("Q20").Value = '$/bandeja 2 kilos'
$ws.Range("R20").Value = 'Provincia de Linares'
$ws.Range("S20").Value = 1250
$ws.Range("T20").Value = 2
